$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 1
$ws.Range("L2").Value = "[0.45454545454545453, 0.2, 0.1, 0.4, 0.2, 0.3, 0.4, 0.1, 0.2, 0.3, 0.1, 0.2, 0.4, 0.2, 0.0, 0.1, 0.2, 0.3, 0.4, 0.5, 0.4, 0.0, 0.4, 0.3, 0.4, 0.6, 0.4, 0.3, 0.3, 0.0, 0.2, 0.1, 0.2, 0.1, 0.3, 0.2, 0.2, 0.0, 0.1, 0.3, 0.2, 0.0, 0.4, 0.2, 0.1, 0.5, 0.2, 0.2, 0.3, 0.5, 0.2, 0.4, 0.1, 0.1, 0.1, 0.1, 0.2, 0.1, 0.1, 0.2, 0.2, 0.1, 0.0, 0.4, 0.3, 0.3, 0.3, 0.1, 0.2, 0.2, 0.4, 0.4, 0.5, 0.2, 0.4, 0.2, 0.2, 0.3, 0.1, 0.3, 0.3, 0.2, 0.1, 0.3, 0.4, 0.3, 0.2, 0.3, 0.0, 0.2, 0.3, 0.2, 0.3, 0.3, 0.2, 0.5, 0.3, 0.2, 0.3, 0.4, 0.3, 0.2, 0.2, 0.1, 0.4, 0.1, 0.0, 0.1, 0.3, 0.2, 0.4, 0.2, 0.0, 0.4, 0.2, 0.2, 0.0, 0.3, 0.4, 0.0, 0.2, 0.2, 0.3, 0.4, 0.2, 0.2, 0.5, 0.4, 0.3, 0.0, 0.3, 0.3, 0.2, 0.4, 0.2, 0.2, 0.5, 0.2, 0.0, 0.3, 0.4, 0.2, 0.1, 0.3, 0.4, 0.1, 0.3, 0.1, 0.3, 0.4, 0.4, 0.2, 0.0, 0.2, 0.3, 0.3, 0.2, 0.6, 0.1, 0.3, 0.3, 0.4, 0.3, 0.3, 0.8, 0.2, 0.4, 0.1, 0.4, 0.2, 0.2, 0.1, 0.2, 0.3, 0.2, 0.3, 0.5, 0.1, 0.6, 0.1, 0.3, 0.1, 0.8, 0.3, 0.4, 0.2, 0.2, 0.4, 0.3, 0.3, 0.4, 0.5, 0.2, 0.2, 0.6, 0.4, 0.4, 0.3, 0.2, 0.0, 0.3, 0.3, 0.4, 0.4, 0.0, 0.4, 0.2, 0.1, 0.2, 0.2, 0.2, 0.2, 0.5, 0.2, 0.5, 0.1, 0.5, 0.5, 0.2, 0.4, 0.4, 0.2, 0.2, 0.2, 0.0, 0.4, 0.3, 0.3, 0.0, 0.2, 0.6, 0.0, 0.3, 0.1, 0.1, 0.1, 0.1, 0.2, 0.2, 0.3, 0.4, 0.2, 0.1, 0.3, 0.2, 0.1, 0.2, 0.4, 0.5, 0.1, 0.1, 0.3, 0.5, 0.1, 0.3, 0.3, 0.1, 0.2, 0.1, 0.2, 0.2, 0.2, 0.0, 0.3, 0.2, 0.0, 0.5, 0.3, 0.2, 0.1, 0.1, 0.0, 0.4, 0.5, 0.4, 0.2, 0.3, 0.0, 0.0, 0.2, 0.1, 0.3, 0.1, 0.4, 0.3, 0.1, 0.0, 0.6, 0.3, 0.3, 0.5, 0.2, 0.2, 0.4, 0.2, 0.2, 0.3, 0.3, 0.2, 0.1, 0.2, 0.3, 0.2, 0.3, 0.1, 0.3, 0.2, 0.2, 0.1, 0.2, 0.1, 0.2, 0.5, 0.1, 0.1, 0.5, 0.3, 0.3, 0.5, 0.3, 0.2, 0.1, 0.1, 0.1, 0.4, 0.3, 0.1, 0.4, 0.5, 0.0, 0.3, 0.2, 0.1, 0.7, 0.1, 0.2, 0.3, 0.3, 0.4, 0.3, 0.3, 0.3, 0.2, 0.2, 0.2, 0.2, 0.0, 0.2, 0.4, 0.3, 0.2, 0.3, 0.1, 0.1, 0.5, 0.3, 0.4, 0.5, 0.3, 0.3, 0.3, 0.2, 0.4, 0.1, 0.2, 0.2, 0.2, 0.3, 0.2, 0.2, 0.2, 0.3, 0.3, 0.1, 0.4, 0.4, 0.1, 0.1, 0.5, 0.2, 0.3, 0.2, 0.2, 0.2, 0.1, 0.1, 0.5, 0.1, 0.1, 0.2, 0.2, 0.1, 0.1, 0.2, 0.3, 0.3, 0.1, 0.3, 0.2, 0.3, 0.1, 0.2, 0.3, 0.3, 0.1, 0.1, 0.1, 0.2, 0.2, 0.1, 0.2, 0.2, 0.1, 0.4, 0.0, 0.1, 0.0, 0.2, 0.1, 0.4, 0.1, 0.1, 0.2, 0.1, 0.5, 0.4, 0.1, 0.3, 0.3, 0.2, 0.4, 0.1, 0.3, 0.3, 0.2, 0.1, 0.4, 0.3, 0.0, 0.5, 0.3, 0.1, 0.1, 0.2, 0.1, 0.2, 0.1, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.3, 0.1, 0.0, 0.2, 0.3, 0.2, 0.2, 0.2, 0.1, 0.4, 0.1, 0.4, 0.3, 0.5, 0.1, 0.5, 0.1, 0.2, 0.2, 0.3, 0.2, 0.2, 0.2, 0.3, 0.4, 0.2, 0.2, 0.3, 0.1, 0.2, 0.3, 0.5, 0.3, 0.1, 0.3, 0.2, 0.5, 0.3, 0.2, 0.3, 0.3, 0.1, 0.4, 0.2, 0.2, 0.1]"
$ws.Range("M2").Value = 0.2600000000000001
$ws.Range("N2").Value = 0.3799999999999999
$ws.Range("O2").Value = 1920

# Update row 3
$ws.Range("A3").Value = 2
$ws.Range("L3").Value = "[0.36363636363636365, 0.4, 0.3, 0.3, 0.3, 0.6, 0.2, 0.2, 0.0, 0.5, 0.3, 0.2, 0.1, 0.2, 0.2, 0.3, 0.4, 0.2, 0.5, 0.1, 0.4, 0.3, 0.1, 0.5, 0.2, 0.4, 0.2, 0.4, 0.4, 0.4, 0.1, 0.3, 0.1, 0.2, 0.2, 0.3, 0.4, 0.3, 0.0, 0.1, 0.2, 0.2, 0.3, 0.1, 0.2, 0.4, 0.1, 0.3, 0.4, 0.1, 0.5, 0.2, 0.2, 0.0, 0.3, 0.2, 0.2, 0.2, 0.3, 0.1, 0.3, 0.3, 0.3, 0.3, 0.4, 0.3, 0.1, 0.1, 0.3, 0.1, 0.5, 0.0, 0.3, 0.2, 0.6, 0.3, 0.4, 0.3, 0.5, 0.2, 0.2, 0.4, 0.5, 0.5, 0.1, 0.2, 0.3, 0.0, 0.1, 0.0, 0.1, 0.2, 0.2, 0.1, 0.6, 0.2, 0.5, 0.4, 0.2, 0.2, 0.4, 0.6, 0.1, 0.2, 0.2, 0.3, 0.4, 0.4, 0.2, 0.3, 0.7, 0.1, 0.2, 0.4, 0.4, 0.3, 0.4, 0.3, 0.2, 0.3, 0.4, 0.2, 0.4, 0.2, 0.3, 0.1, 0.2, 0.0, 0.3, 0.5, 0.4, 0.3, 0.3, 0.4, 0.4, 0.1, 0.2, 0.2, 0.2, 0.5, 0.5, 0.5, 0.3, 0.5, 0.4, 0.1, 0.2, 0.4, 0.3, 0.4, 0.1, 0.4, 0.3, 0.2, 0.2, 0.3, 0.3, 0.2, 0.3, 0.6, 0.3, 0.2, 0.3, 0.2, 0.2, 0.2, 0.1, 0.3, 0.5, 0.4, 0.2, 0.3, 0.4, 0.2, 0.3, 0.4, 0.4, 0.4, 0.0, 0.2, 0.2, 0.2, 0.3, 0.3, 0.4, 0.2, 0.5, 0.2, 0.2, 0.1, 0.1, 0.0, 0.3, 0.1, 0.4, 0.3, 0.5, 0.2, 0.1, 0.4, 0.4, 0.2, 0.3, 0.3, 0.1, 0.3, 0.1, 0.3, 0.4, 0.3, 0.5, 0.2, 0.5, 0.7, 0.1, 0.1, 0.1, 0.4, 0.2, 0.4, 0.4, 0.2, 0.3, 0.0, 0.3, 0.2, 0.4, 0.3, 0.3, 0.5, 0.2, 0.5, 0.3, 0.3, 0.2, 0.3, 0.4, 0.1, 0.5, 0.2, 0.3, 0.4, 0.2, 0.1, 0.0, 0.3, 0.2, 0.4, 0.4, 0.2, 0.2, 0.4, 0.4, 0.5, 0.3, 0.3, 0.2, 0.3, 0.3, 0.2, 0.3, 0.2, 0.4, 0.2, 0.3, 0.3, 0.4, 0.1, 0.2, 0.1, 0.1, 0.2, 0.5, 0.3, 0.1, 0.1, 0.3, 0.2, 0.1, 0.1, 0.1, 0.2, 0.3, 0.2, 0.2, 0.2, 0.2, 0.5, 0.5, 0.1, 0.1, 0.2, 0.5, 0.3, 0.3, 0.2, 0.4, 0.4, 0.1, 0.0, 0.1, 0.5, 0.2, 0.3, 0.1, 0.4, 0.0, 0.0, 0.3, 0.1, 0.1, 0.2, 0.1, 0.2, 0.3, 0.1, 0.2, 0.4, 0.4, 0.1, 0.2, 0.1, 0.7, 0.0, 0.3, 0.4, 0.1, 0.1, 0.4, 0.2, 0.1, 0.3, 0.3, 0.3, 0.2, 0.2, 0.3, 0.4, 0.4, 0.2, 0.3, 0.2, 0.2, 0.1, 0.3, 0.2, 0.4, 0.4, 0.1, 0.4, 0.1, 0.4, 0.2, 0.1, 0.1, 0.3, 0.3, 0.2, 0.4, 0.2, 0.2, 0.5, 0.1, 0.2, 0.2, 0.1, 0.1, 0.4, 0.4, 0.5, 0.0, 0.3, 0.2, 0.3, 0.2, 0.2, 0.4, 0.1, 0.3, 0.4, 0.3, 0.3, 0.2, 0.1, 0.2, 0.2, 0.3, 0.2, 0.4, 0.1, 0.3, 0.2, 0.4, 0.4, 0.3, 0.2, 0.1, 0.5, 0.2, 0.1, 0.3, 0.4, 0.1, 0.1, 0.2, 0.2, 0.2, 0.2, 0.1, 0.4, 0.3, 0.0, 0.1, 0.5, 0.2, 0.3, 0.5, 0.3, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.3, 0.2, 0.3, 0.6, 0.2, 0.2, 0.1, 0.6, 0.1, 0.3, 0.2, 0.4, 0.5, 0.4, 0.3, 0.4, 0.2, 0.4, 0.2, 0.2, 0.4, 0.4, 0.2, 0.3, 0.2, 0.2, 0.4, 0.3, 0.2, 0.7, 0.4, 0.4, 0.1, 0.3, 0.2, 0.3, 0.2, 0.2, 0.2, 0.6, 0.2, 0.1, 0.1, 0.2, 0.2, 0.2, 0.1, 0.5, 0.1, 0.2, 0.4, 0.3, 0.3, 0.3, 0.2, 0.4, 0.6, 0.1, 0.3, 0.1, 0.1, 0.1, 0.2, 0.2, 0.3, 0.4, 0.4, 0.3, 0.1, 0.2, 0.4, 0.0, 0.1, 0.3, 0.0]"
$ws.Range("M3").Value = 0.22
$ws.Range("N3").Value = 0.39
$ws.Range("O3").Value = 840

# Delete rows 4-8 (old data no longer present)
$ws.Range("A4:P8").EntireRow.Delete()
